$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# --- Step 1: free up shared-string slot 41 (currently "New 200m transport / depth",
#     sole referrer is A15) so the new strings we add land on the same indices
#     the target workbook uses. ---
$ws.Range("A15").ClearContents()

# --- Step 2: add the four new rows (22-25), copying number formats from the
#     closest existing template rows so style indices match. ---

# Row 22 <- format like row 21 (plain numeric style throughout)
$ws.Range("A21:I21").Copy()
$ws.Range("A22:I22").PasteSpecial(-4122)
$ws.Range("A22:I22").RowHeight = 32
$ws.Range("A22").Value = "New 200m transport / depth, dt = 1 hr, j = 2, swim to shallow"
$ws.Range("C22").Value = 0.33189999999999997
$ws.Range("E22").Value = 15.1043
$ws.Range("I22").Value = 2.1004

# Row 23 <- format like row 20 (scientific-notation style on C, E, I)
$ws.Range("A20:I20").Copy()
$ws.Range("A23:I23").PasteSpecial(-4122)
$ws.Range("A23:I23").RowHeight = 32
$ws.Range("A23").Value = "New 200m transport / depth, dt = 1 hr, j = 2, swim to deep"
$ws.Range("C23").Value = 7470400000
$ws.Range("E23").Value = [double]"6.6192999999999997E+46"
$ws.Range("I23").Value = [double]"5.4772000000000002E+45"

# --- Step 3: restore A15 with its new text - lands on the freed slot. ---
$ws.Range("A15").Value = "New 200m transport / depth, dt = 1 hr, j = 2"

# Row 24 <- format like row 21, but I24 needs scientific-notation style (like I20)
$ws.Range("A21:I21").Copy()
$ws.Range("A24:I24").PasteSpecial(-4122)
$ws.Range("A24:I24").RowHeight = 32
$ws.Range("I20").Copy()
$ws.Range("I24").PasteSpecial(-4122)
$ws.Range("A24").Value = "New 200m transport / depth, dt = 1 hr, j = 2, swim to const rand"
$ws.Range("I24").Value = [double]"2.2199000000000001E+46"

# Row 25 <- format like row 21
$ws.Range("A21:I21").Copy()
$ws.Range("A25:I25").PasteSpecial(-4122)
$ws.Range("A25:I25").RowHeight = 32
$ws.Range("A25").Value = "New 200m transport / depth, dt = 1 hr, j = 2, swim to changing rand"
$ws.Range("I25").Value = 44.083599999999997

$excel.CutCopyMode = 0

# --- Step 4: selection / view state ---
$ws.Range("I26").Select() | Out-Null
$win = $excel.ActiveWindow
$win.ScrollRow = 5
$win.ScrollColumn = 1
$win.ScrollColumn | Out-Null
